$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the ANTONIO (004241147) row - it is row 3 (row 1 = header, row 2 = MARCUS)
$ws.Rows.Item(3).Delete()

# 2. Update DIEGO's (004479965) balance from 11173.5 to 16173.58
$ws.Range("C7").Value = 16173.58

# 3. Insert a new row for MARIANA (005000460 / 3299.99) right before PEDRO (005232019),
#    which currently sits at row 9 after the deletion above.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "005000460"
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = "MARIANA"
$ws.Range("C9").Value = 3299.99

# 4. Insert a new row for OTAVIO (004452946 / 2000) right after CARLA (004643153),
#    which is now at row 11.
$ws.Rows.Item(12).Insert()
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "004452946"
$ws.Range("A12").ClearFormats()
$ws.Range("B12").Value = "OTAVIO"
$ws.Range("C12").Value = 2000
